$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Qty executed upto date" (column C) values
$ws.Range("C8").Value = 18
$ws.Range("C9").Value = 51
$ws.Range("C10").Value = 79
$ws.Range("C11").Value = 58
$ws.Range("C12").Value = 7
$ws.Range("C13").Value = 81
$ws.Range("C14").Value = 8
$ws.Range("C15").Value = 76
$ws.Range("C16").Value = 73
$ws.Range("C17").Value = 61

# Update corresponding "Upto date Amount" (column G) string values (Rate * Qty, formatted to 2 decimals)
# These cells hold their amount as text (e.g. "13056.00"), so force Text format
# first to stop Excel from auto-converting the numeric-looking string to a number.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "13056.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "37288.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "38396.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "11016.00"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "184.00"

# Update Grand Total rows
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "99940.00"

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "99940.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "99940.00"

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "99940.00"
